$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.655.86"
$ws.Range("E2").Value = "  -2.40%  "
$ws.Range("D3").Value = "2.336.30"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'500.45"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "'128.62"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -3.28%  "
$ws.Range("D9").Value = "2.342.26"
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("D10").Value = "'0.0979"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "'4.83"
$ws.Range("E12").Value = "  +4.67%  "
$ws.Range("D13").Value = "'0.321"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "2.748.94"
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "55.639.32"
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("D16").Value = "'21.54"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").Value = "2.322.20"
$ws.Range("E18").Value = "  -5.03%  "
$ws.Range("D19").Value = "'9.89"
$ws.Range("E19").Value = "  -3.81%  "
$ws.Range("D20").Value = "'308.16"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "'4.00"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("D22").Value = "'6.19"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'65.24"
$ws.Range("E24").Value = "  -3.82%  "
$ws.Range("D25").Value = "'0.992"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("D26").Value = "'0.371"
$ws.Range("E26").Value = "  -1.50%  "
$ws.Range("D27").Value = "'0.146"
$ws.Range("E27").Value = "  -4.13%  "
$ws.Range("D28").Value = "'7.13"
$ws.Range("E28").Value = "  -4.18%  "
$ws.Range("D29").Value = "'172.79"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").Value = "'1.63"
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("D31").Value = "0.0₃0698"
$ws.Range("E31").Value = "  -3.95%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "'5.79"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  -5.78%  "
$ws.Range("D36").Value = "'17.57"
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").Value = "'0.825"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'3.64"
$ws.Range("E39").Value = "  -5.20%  "
$ws.Range("D40").Value = "'36.09"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").Value = "'1.38"
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").Value = "'127.11"
$ws.Range("E43").Value = "  -4.21%  "
$ws.Range("D44").Value = "'4.69"
$ws.Range("E44").Value = "  -5.78%  "
$ws.Range("D45").Value = "'0.555"
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("D46").Value = "'0.0892"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").Value = "'236.15"
$ws.Range("E47").Value = "  -6.03%  "
$ws.Range("D48").Value = "'0.0477"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").Value = "'0.0205"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").Value = "'16.69"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("E51").Value = "  -0.24%  "
